$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Header date
Replace-Text "2025-05-25 Sunday" "2025-05-26 Monday"

# Table row 1
Replace-Text "23÷7=" "57÷5="
Replace-Text "42÷7=" "92÷9="
Replace-Text "22÷6=" "71÷8="
Replace-Text "10÷7=" "95÷3="
Replace-Text "63÷9=" "98÷7="

# Table row 5
Replace-Text "88÷2=" "98÷2="
Replace-Text "60÷8=" "72÷3="
Replace-Text "28÷2=" "67÷6="
Replace-Text "13÷2=" "59÷7="
Replace-Text "74÷5=" "41÷4="

# Table row 9
Replace-Text "26÷2=" "66÷5="
Replace-Text "34÷4=" "61÷8="
Replace-Text "41÷3=" "98÷8="
Replace-Text "94÷7=" "16÷9="
Replace-Text "27÷5=" "42÷6="

# Table row 13 (first cell has duplicate text "19÷2=" elsewhere in doc,
# so address it directly via the table cell rather than Find/Replace)
$tbl = $d.Tables.Item(1)
$tbl.Cell(13, 1).Range.Text = "88÷4="
Replace-Text "89÷8=" "31÷4="
Replace-Text "58÷7=" "59÷5="
Replace-Text "70÷9=" "91÷4="
Replace-Text "85÷2=" "25÷2="

# Table row 17
Replace-Text "99÷3=" "29÷9="
Replace-Text "56÷3=" "63÷3="
$tbl.Cell(17, 3).Range.Text = "66÷4="
Replace-Text "43÷8=" "76÷3="
Replace-Text "36÷2=" "63÷8="
